# Update "想去人数" (F column) counts across sheets, as scraped at a later run.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 301
$ws1.Range("F4").Value  = 1235
$ws1.Range("F5").Value  = 351
$ws1.Range("F6").Value  = 315
$ws1.Range("F7").Value  = 3818
$ws1.Range("F10").Value = 1778
$ws1.Range("F11").Value = 326
$ws1.Range("F13").Value = 727
$ws1.Range("F14").Value = 147
$ws1.Range("F16").Value = 2069
$ws1.Range("F20").Value = 324
$ws1.Range("F23").Value = 266

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value  = 37
$ws2.Range("F9").Value  = 125
$ws2.Range("F10").Value = 89

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F9").Value  = 37
$ws4.Range("F10").Value = 37
$ws4.Range("F12").Value = 301
$ws4.Range("F13").Value = 1236
$ws4.Range("F14").Value = 351
$ws4.Range("F18").Value = 315
$ws4.Range("F19").Value = 3818
$ws4.Range("F20").Value = 125
$ws4.Range("F22").Value = 89
$ws4.Range("F26").Value = 1778
$ws4.Range("F27").Value = 326
$ws4.Range("F30").Value = 727
$ws4.Range("F31").Value = 147
$ws4.Range("F34").Value = 2069
$ws4.Range("F40").Value = 324
$ws4.Range("F50").Value = 266
